$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update buyIn values in row 2
$ws.Range("B2").Value = "100-5;rated-6009;1"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""

# Update wager / player values in row 2
$ws.Range("I2").Value = "P1;300;P2;TGT"
$ws.Range("J2").Value = "P4;100;B5"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""

# Update TakeBets / payAmt values in row 2
$ws.Range("V2").Value = "B5;TGT"
$ws.Range("W2").Value = ""
$ws.Range("X2").Value = ""
$ws.Range("Y2").Value = ""

# Update selection to match target state
$ws.Range("Y8").Select()
